$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Query text in column C for the rows whose text changed.
$ws.Range("C3").Value  = "Who sings 'Photograph'?"
$ws.Range("C6").Value  = 'Find 2 songs that have a color in the lyrics and the artist who sang them'
$ws.Range("C7").Value  = 'What is the song that has similar lyrics to "Never Grow Up" by Taylor Swift?'
$ws.Range("C9").Value  = "How is love mentioned in songs?"
$ws.Range("C11").Value = 'Who was mentioned in song "Dear John" by Taylor Swift?'
$ws.Range("C13").Value = 'Write a poem that is inspired by song "Thinking Out Loud" by Ed Sheeran'

# Widen column C slightly and add a custom width for (currently empty) column D.
$ws.Columns.Item(3).ColumnWidth = 59.666666666666664
$ws.Columns.Item(4).ColumnWidth = 57.5

# Move the active selection to C10, as recorded in the saved view state.
$ws.Range("C10").Select()
